$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 256.25
$ws.Range("I28").Value = 256.25
$ws.Range("K28").Value = 256.25
$ws.Range("M28").Value = 228.75
$ws.Range("H31").Value = 1089.9
$ws.Range("I31").Value = 1089.9
$ws.Range("K31").Value = 3269.7
$ws.Range("M31").Value = -3039.7
$ws.Range("H43").Value = 2500
$ws.Range("J43").Value = 2500
$ws.Range("L43").Value = 2500
$ws.Range("N43").Value = -2638
$ws.Range("H100").Value = 1449.75
$ws.Range("I100").Value = 1449.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1449.75
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -908.75
$ws.Range("N100").ClearContents()
$ws.Range("H103").Value = 198.5
$ws.Range("J103").Value = 198.5
$ws.Range("L103").Value = 595.5
$ws.Range("N103").Value = -1767.5
$ws.Range("H104").Value = 169.5
$ws.Range("I104").Value = 169.5
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 508.5
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 1238.5
$ws.Range("N104").ClearContents()
$ws.Range("H106").Value = 38336.5
$ws.Range("I106").Value = 50745
$ws.Range("K106").Value = 50745
$ws.Range("M106").Value = -50114
$ws.Range("H107").Value = 244.28572
$ws.Range("I107").Value = 142.8
$ws.Range("K107").Value = 142.8
$ws.Range("M107").Value = 1777.2
$ws.Range("H111").Value = 47141.855
$ws.Range("I111").Value = 46399.8
$ws.Range("J111").Value = 48997
$ws.Range("K111").Value = 139199.4
$ws.Range("L111").Value = 146991
$ws.Range("M111").Value = -136132.4
$ws.Range("N111").Value = -153125
$ws.Range("I113").Value = 4498.5
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 4498.5
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = -1244.5
$ws.Range("N113").Value = -11008
$ws.Range("H132").Value = 2172.2
$ws.Range("J132").Value = 4990
$ws.Range("L132").Value = 14970
$ws.Range("N132").Value = -20030
$ws.Range("H138").Value = 3673.1052
$ws.Range("I138").Value = 3287.25
$ws.Range("K138").Value = 9861.75
$ws.Range("M138").Value = -4721.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H45").Value = 735
$ws.Range("I45").Value = 735
$ws.Range("K45").Value = 735
$ws.Range("M45").Value = -358
$ws.Range("H61").Value = 2499.889
$ws.Range("I61").Value = 2187.5
$ws.Range("K61").Value = 2187.5
$ws.Range("M61").Value = -1975.5
$ws.Range("H74").Value = 5124.5
$ws.Range("I74").Value = 7999
$ws.Range("J74").Value = 2250
$ws.Range("K74").Value = 7999
$ws.Range("L74").Value = 2250
$ws.Range("M74").Value = -7125
$ws.Range("N74").Value = -3998
$ws.Range("H77").Value = 5124.5
$ws.Range("I77").Value = 7999
$ws.Range("J77").Value = 2250
$ws.Range("K77").Value = 39995
$ws.Range("L77").Value = 11250
$ws.Range("M77").Value = -35627
$ws.Range("N77").Value = -19986
$ws.Range("H97").Value = 650
$ws.Range("I97").Value = 660
$ws.Range("K97").Value = 660
$ws.Range("M97").Value = -164
$ws.Range("H110").Value = 6167625
$ws.Range("I110").Value = 6167625
$ws.Range("K110").Value = 6167625
$ws.Range("M110").Value = -6165580
$ws.Range("H132").Value = 2892.2
$ws.Range("I132").Value = 2892.2
$ws.Range("K132").Value = 8676.599999999999
$ws.Range("M132").Value = -6146.599999999999
$ws.Range("H136").Value = 2499.889
$ws.Range("I136").Value = 2187.5
$ws.Range("K136").Value = 6562.5
$ws.Range("M136").Value = -4012.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 3
$ws.Range("J11").Value = 3
$ws.Range("L11").Value = 3
$ws.Range("N11").Value = -283
$ws.Range("H26").Value = 4494.5
$ws.Range("I26").Value = 4494.5
$ws.Range("K26").Value = 4494.5
$ws.Range("M26").Value = -4202.5
$ws.Range("H134").Value = 2100.3572
$ws.Range("I134").Value = 2030.0769
$ws.Range("K134").Value = 6090.2307
$ws.Range("M134").Value = -3555.2307

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2375.5557
$ws.Range("I31").Value = 2172.5
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 2172.5
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -1877.5
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 2375.5557
$ws.Range("I34").Value = 2172.5
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 2172.5
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1970.5
$ws.Range("N34").Value = -4404
$ws.Range("H41").Value = 900
$ws.Range("I41").Value = 900
$ws.Range("K41").Value = 900
$ws.Range("M41").Value = -472
$ws.Range("H58").Value = 1786.3636
$ws.Range("I58").Value = 1589.2222
$ws.Range("J58").Value = 2673.5
$ws.Range("K58").Value = 1589.2222
$ws.Range("L58").Value = 2673.5
$ws.Range("M58").Value = -1386.2222
$ws.Range("N58").Value = -3079.5
$ws.Range("H99").Value = 1444.9412
$ws.Range("I99").Value = 1314
$ws.Range("K99").Value = 1314
$ws.Range("M99").Value = 184
$ws.Range("H126").Value = 1444.9412
$ws.Range("I126").Value = 1314
$ws.Range("K126").Value = 3942
$ws.Range("M126").Value = -1472
$ws.Range("H132").Value = 9326.333000000001
$ws.Range("I132").Value = 9326.333000000001
$ws.Range("K132").Value = 27978.999
$ws.Range("M132").Value = -25448.999
$ws.Range("H134").Value = 3600.8667
$ws.Range("I134").Value = 2978
$ws.Range("J134").Value = 4535.1665
$ws.Range("K134").Value = 8934
$ws.Range("L134").Value = 13605.4995
$ws.Range("M134").Value = -6399
$ws.Range("N134").Value = -18675.4995
$ws.Range("H136").Value = 1786.3636
$ws.Range("I136").Value = 1589.2222
$ws.Range("J136").Value = 2673.5
$ws.Range("K136").Value = 4767.6666
$ws.Range("L136").Value = 8020.5
$ws.Range("M136").Value = -2217.6666
$ws.Range("N136").Value = -13120.5
$ws.Range("H141").Value = 36241.215
$ws.Range("J141").Value = 36241.215
$ws.Range("L141").Value = 36241.215
$ws.Range("N141").Value = -46601.215

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1159.5
$ws.Range("J23").Value = 1159.5
$ws.Range("L23").Value = 3478.5
$ws.Range("N23").Value = -3948.5
$ws.Range("H45").Value = 2033
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 51623.668
$ws.Range("J15").Value = 51623.668
$ws.Range("L15").Value = 51623.668
$ws.Range("N15").Value = -52199.668
$ws.Range("H81").Value = 51623.668
$ws.Range("J81").Value = 51623.668
$ws.Range("L81").Value = 51623.668
$ws.Range("N81").Value = -53619.668
$ws.Range("H84").Value = 51623.668
$ws.Range("J84").Value = 51623.668
$ws.Range("L84").Value = 154871.004
$ws.Range("N84").Value = -164855.004
$ws.Range("H99").Value = 8299.571
$ws.Range("I99").Value = 5538.25
$ws.Range("K99").Value = 5538.25
$ws.Range("M99").Value = -3292.25
$ws.Range("H113").Value = 699
$ws.Range("H122").Value = 8749.666999999999
$ws.Range("I122").Value = 7666.6665
$ws.Range("J122").Value = 9832.666999999999
$ws.Range("K122").Value = 22999.9995
$ws.Range("L122").Value = 29498.001
$ws.Range("M122").Value = -20549.9995
$ws.Range("N122").Value = -34398.001
$ws.Range("H126").Value = 2666.6667
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 7736
$ws.Range("I132").Value = 8089
$ws.Range("J132").Value = 2794
$ws.Range("K132").Value = 24267
$ws.Range("L132").Value = 8382
$ws.Range("M132").Value = -21737
$ws.Range("N132").Value = -13442

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7849.278
$ws.Range("I7").Value = 9058.799999999999
$ws.Range("K7").Value = 9058.799999999999
$ws.Range("M7").Value = -8946.799999999999
$ws.Range("H126").Value = 7849.278
$ws.Range("I126").Value = 9058.799999999999
$ws.Range("K126").Value = 27176.4
$ws.Range("M126").Value = -24706.4
$ws.Range("H132").Value = 2062.6
$ws.Range("I132").Value = 2402
$ws.Range("J132").Value = 705
$ws.Range("K132").Value = 7206
$ws.Range("L132").Value = 2115
$ws.Range("M132").Value = -4676
$ws.Range("N132").Value = -7175
$ws.Range("H136").Value = 3417.875
$ws.Range("I136").Value = 3417.875
$ws.Range("K136").Value = 10253.625
$ws.Range("M136").Value = -7703.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H41").Value = 17973.625
$ws.Range("I41").Value = 16658
$ws.Range("K41").Value = 16658
$ws.Range("M41").Value = -16268
$ws.Range("H49").Value = 50000
$ws.Range("J49").Value = 50000
$ws.Range("L49").Value = 50000
$ws.Range("N49").Value = -50460
$ws.Range("H119").Value = 84999.5
$ws.Range("J119").Value = 84999.5
$ws.Range("L119").Value = 84999.5
$ws.Range("N119").Value = -94675.5
$ws.Range("H132").Value = 1504
$ws.Range("I132").Value = 1504
$ws.Range("K132").Value = 4512
$ws.Range("M132").Value = -1982
$ws.Range("H136").Value = 1895.0588
$ws.Range("I136").Value = 1951
$ws.Range("K136").Value = 5853
$ws.Range("M136").Value = -3303

Write-Host "Applied all changes"